$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17
$ws.Cells.Item(17, 1).Value = 'PS Number'
$ws.Cells.Item(17, 2).Value = 'Mathematics'
$ws.Cells.Item(17, 3).Value = 'Hindi'
$ws.Cells.Item(17, 4).Value = 'Biology'
$ws.Cells.Item(17, 5).Value = 'Physics'
$ws.Cells.Item(17, 6).Value = 'Telugu'
$ws.Cells.Item(17, 7).Value = 'English'
$ws.Cells.Item(17, 8).Value = 'Chemistry'
$ws.Cells.Item(17, 9).Value = 'C-Programming'
$ws.Cells.Item(17, 10).Value = 'Java'
$ws.Cells.Item(17, 11).Value = 'Processor and Controllers'
$ws.Cells.Item(17, 12).Value = 'Analog Communication'
$ws.Cells.Item(17, 13).Value = 'Digital Communication'
$ws.Cells.Item(17, 14).Value = 'Operating systems'
$ws.Cells.Item(17, 15).Value = 'Python'
$ws.Cells.Item(17, 16).Value = 'Web Technologies'
$ws.Cells.Item(17, 17).Value = 'Engineering Drawing'
$ws.Cells.Item(17, 18).Value = 'Geography'
$ws.Cells.Item(17, 19).Value = 'History'
$ws.Cells.Item(17, 20).Value = 'Civics '
$ws.Cells.Item(17, 21).Value = 'Economics'

# Row 18
$ws.Cells.Item(18, 1).Value = 99004400
$ws.Cells.Item(18, 2).Value = 95
$ws.Cells.Item(18, 3).Value = 98
$ws.Cells.Item(18, 4).Value = 87
$ws.Cells.Item(18, 5).Value = 79
$ws.Cells.Item(18, 6).Value = 65
$ws.Cells.Item(18, 7).Value = 90
$ws.Cells.Item(18, 8).Value = 87
$ws.Cells.Item(18, 9).Value = 95
$ws.Cells.Item(18, 10).Value = 74
$ws.Cells.Item(18, 11).Value = 86
$ws.Cells.Item(18, 12).Value = 75
$ws.Cells.Item(18, 13).Value = 88
$ws.Cells.Item(18, 14).Value = 76
$ws.Cells.Item(18, 15).Value = 87
$ws.Cells.Item(18, 16).Value = 66
$ws.Cells.Item(18, 17).Value = 79
$ws.Cells.Item(18, 18).Value = 80
$ws.Cells.Item(18, 19).Value = 82
$ws.Cells.Item(18, 20).Value = 76
$ws.Cells.Item(18, 21).Value = 87

# Row 19
$ws.Cells.Item(19, 1).Value = 'PS Number'
$ws.Cells.Item(19, 2).Value = 'Hobby 1'
$ws.Cells.Item(19, 3).Value = 'Hobby 2'
$ws.Cells.Item(19, 4).Value = 'Hobby 3'
$ws.Cells.Item(19, 5).Value = 'Hobby 4'
$ws.Cells.Item(19, 6).Value = 'Hobby 5'
$ws.Cells.Item(19, 7).Value = 'Hobby 6'
$ws.Cells.Item(19, 8).Value = 'Hobby 7'
$ws.Cells.Item(19, 9).Value = 'Hobby 8'
$ws.Cells.Item(19, 10).Value = 'Hobby 9'
$ws.Cells.Item(19, 11).Value = 'Hobby 10'
$ws.Cells.Item(19, 12).Value = 'Hobby 11'
$ws.Cells.Item(19, 13).Value = 'Hobby 12'
$ws.Cells.Item(19, 14).Value = 'Hobby 13'
$ws.Cells.Item(19, 15).Value = 'Hobby 14'
$ws.Cells.Item(19, 16).Value = 'Hobby 15'
$ws.Cells.Item(19, 17).Value = 'Hobby 16'
$ws.Cells.Item(19, 18).Value = 'Hobby 17'
$ws.Cells.Item(19, 19).Value = 'Hobby 18'
$ws.Cells.Item(19, 20).Value = 'Hobby 19'
$ws.Cells.Item(19, 21).Value = 'Hobby 20'

# Row 20
$ws.Cells.Item(20, 1).Value = 99004408
$ws.Cells.Item(20, 2).Value = 'Go Camping'
$ws.Cells.Item(20, 3).Value = 'Watch Documentaries'
$ws.Cells.Item(20, 4).Value = 'New Music Discovery'
$ws.Cells.Item(20, 5).Value = 'Computer Programming'
$ws.Cells.Item(20, 6).Value = 'Fishing'
$ws.Cells.Item(20, 7).Value = 'Walking'
$ws.Cells.Item(20, 8).Value = 'Travelling'
$ws.Cells.Item(20, 9).Value = 'Golf'
$ws.Cells.Item(20, 10).Value = 'Exercise'
$ws.Cells.Item(20, 11).Value = 'Drawing'
$ws.Cells.Item(20, 12).Value = 'Sewing'
$ws.Cells.Item(20, 13).Value = 'Hiking'
$ws.Cells.Item(20, 14).Value = 'Cooking'
$ws.Cells.Item(20, 15).Value = 'Scrapbooking'
$ws.Cells.Item(20, 16).Value = 'Cross-Stitch'
$ws.Cells.Item(20, 17).Value = 'Jigsaw puzzel'
$ws.Cells.Item(20, 18).Value = 'Parachuting'
$ws.Cells.Item(20, 19).Value = 'Marketing'
$ws.Cells.Item(20, 20).Value = 'Genealogy'
$ws.Cells.Item(20, 21).Value = 'Computer Programming'

# Row 21
$ws.Cells.Item(21, 1).Value = 'PS Number'
$ws.Cells.Item(21, 2).Value = 'City 1'
$ws.Cells.Item(21, 3).Value = ' City 2'
$ws.Cells.Item(21, 4).Value = ' City 3'
$ws.Cells.Item(21, 5).Value = 'City 4'
$ws.Cells.Item(21, 6).Value = 'City 5'
$ws.Cells.Item(21, 7).Value = 'City 6'
$ws.Cells.Item(21, 8).Value = 'City 7'
$ws.Cells.Item(21, 9).Value = 'City 8'
$ws.Cells.Item(21, 10).Value = 'City 9'
$ws.Cells.Item(21, 11).Value = 'City 10'
$ws.Cells.Item(21, 12).Value = 'City 11'
$ws.Cells.Item(21, 13).Value = 'City 12'
$ws.Cells.Item(21, 14).Value = 'City 13'
$ws.Cells.Item(21, 15).Value = 'City 14'
$ws.Cells.Item(21, 16).Value = 'City 15'
$ws.Cells.Item(21, 17).Value = 'City 16'
$ws.Cells.Item(21, 18).Value = 'City 17'
$ws.Cells.Item(21, 19).Value = 'City 18'
$ws.Cells.Item(21, 20).Value = 'City 19'
$ws.Cells.Item(21, 21).Value = 'City 20'

# Row 22
$ws.Cells.Item(22, 1).Value = 99004412
$ws.Cells.Item(22, 2).Value = 'Gwalior'
$ws.Cells.Item(22, 3).Value = 'Puducherry'
$ws.Cells.Item(22, 4).Value = 'Noida'
$ws.Cells.Item(22, 5).Value = 'Jamshedpur'
$ws.Cells.Item(22, 6).Value = 'Cuttak'
$ws.Cells.Item(22, 7).Value = 'Ajmer'
$ws.Cells.Item(22, 8).Value = 'Jammu'
$ws.Cells.Item(22, 9).Value = 'Tirupathi'
$ws.Cells.Item(22, 10).Value = 'Kakinada'
$ws.Cells.Item(22, 11).Value = 'Khammam'
$ws.Cells.Item(22, 12).Value = 'Nizamabad'
$ws.Cells.Item(22, 13).Value = 'Thiruvananthapur'
$ws.Cells.Item(22, 14).Value = 'Guwahati'
$ws.Cells.Item(22, 15).Value = 'Meerut'
$ws.Cells.Item(22, 16).Value = 'Raipur'
$ws.Cells.Item(22, 17).Value = 'Shimla'
$ws.Cells.Item(22, 18).Value = 'Jodhpur'
$ws.Cells.Item(22, 19).Value = 'Ujjain'
$ws.Cells.Item(22, 20).Value = 'Udhaipur'
$ws.Cells.Item(22, 21).Value = 'Gandhinagar'

# Row 23
$ws.Cells.Item(23, 1).Value = 'PS Number'
$ws.Cells.Item(23, 2).Value = 'C'
$ws.Cells.Item(23, 3).Value = 'Python'
$ws.Cells.Item(23, 4).Value = 'Java'
$ws.Cells.Item(23, 5).Value = 'Java Script'
$ws.Cells.Item(23, 6).Value = 'Go'
$ws.Cells.Item(23, 7).Value = 'Perl'
$ws.Cells.Item(23, 8).Value = 'Ruby'
$ws.Cells.Item(23, 9).Value = 'Swift'
$ws.Cells.Item(23, 10).Value = 'Scala'
$ws.Cells.Item(23, 11).Value = 'PHP'
$ws.Cells.Item(23, 12).Value = 'C++'
$ws.Cells.Item(23, 13).Value = 'R'
$ws.Cells.Item(23, 14).Value = 'SQL'
$ws.Cells.Item(23, 15).Value = 'Arduino'
$ws.Cells.Item(23, 16).Value = 'Matlab'
$ws.Cells.Item(23, 17).Value = 'Rust'
$ws.Cells.Item(23, 18).Value = 'Type Script'
$ws.Cells.Item(23, 19).Value = 'Kotlin'
$ws.Cells.Item(23, 20).Value = 'CSS'
$ws.Cells.Item(23, 21).Value = 'Powershell'

# Row 24
$ws.Cells.Item(24, 1).Value = 99004402
$ws.Cells.Item(24, 2).Value = 'Beginner'
$ws.Cells.Item(24, 3).Value = 'Competent'
$ws.Cells.Item(24, 4).Value = 'Beginner'
$ws.Cells.Item(24, 5).Value = 'Expert'
$ws.Cells.Item(24, 6).Value = 'Novice'
$ws.Cells.Item(24, 7).Value = 'Beginner'
$ws.Cells.Item(24, 8).Value = 'Competent'
$ws.Cells.Item(24, 9).Value = 'Competent'
$ws.Cells.Item(24, 10).Value = 'Beginner'
$ws.Cells.Item(24, 11).Value = 'Competent'
$ws.Cells.Item(24, 12).Value = 'Beginner'
$ws.Cells.Item(24, 13).Value = 'Beginner'
$ws.Cells.Item(24, 14).Value = 'Novice'
$ws.Cells.Item(24, 15).Value = 'Beginner'
$ws.Cells.Item(24, 16).Value = 'Proficient'
$ws.Cells.Item(24, 17).Value = 'Beginner'
$ws.Cells.Item(24, 18).Value = 'Novice'
$ws.Cells.Item(24, 19).Value = 'Beginner'
$ws.Cells.Item(24, 20).Value = 'Beginner'
$ws.Cells.Item(24, 21).Value = 'Beginner'

# Row 25
$ws.Cells.Item(25, 1).Value = 'PS Number'
$ws.Cells.Item(25, 2).Value = 'Sport 1'
$ws.Cells.Item(25, 3).Value = 'Sport 2'
$ws.Cells.Item(25, 4).Value = 'Sport 3'
$ws.Cells.Item(25, 5).Value = 'Sport 4'
$ws.Cells.Item(25, 6).Value = 'Sport 5'
$ws.Cells.Item(25, 7).Value = 'Sport 6'
$ws.Cells.Item(25, 8).Value = 'Sport 7'
$ws.Cells.Item(25, 9).Value = 'Sport 8'
$ws.Cells.Item(25, 10).Value = 'Sport 9 '
$ws.Cells.Item(25, 11).Value = 'Sport 10'
$ws.Cells.Item(25, 12).Value = 'Sport 11'
$ws.Cells.Item(25, 13).Value = 'Sport 12'
$ws.Cells.Item(25, 14).Value = 'Sport 13'
$ws.Cells.Item(25, 15).Value = 'Sport 14'
$ws.Cells.Item(25, 16).Value = 'Sport 15'
$ws.Cells.Item(25, 17).Value = 'Sport 16'
$ws.Cells.Item(25, 18).Value = 'Sport 17'
$ws.Cells.Item(25, 19).Value = 'Sport 18'
$ws.Cells.Item(25, 20).Value = 'Sport 19'
$ws.Cells.Item(25, 21).Value = 'Sport 20'

# Row 26
$ws.Cells.Item(26, 1).Value = 99004402
$ws.Cells.Item(26, 2).Value = 'Bungee Jumping'
$ws.Cells.Item(26, 3).Value = 'Hockey'
$ws.Cells.Item(26, 4).Value = 'Yoga'
$ws.Cells.Item(26, 5).Value = 'Car Racing'
$ws.Cells.Item(26, 6).Value = 'Cricket'
$ws.Cells.Item(26, 7).Value = 'Chess'
$ws.Cells.Item(26, 8).Value = 'Carroms'
$ws.Cells.Item(26, 9).Value = 'Dodgeball'
$ws.Cells.Item(26, 10).Value = 'Darts'
$ws.Cells.Item(26, 11).Value = 'Handball'
$ws.Cells.Item(26, 12).Value = 'Hockey'
$ws.Cells.Item(26, 13).Value = 'Kendo'
$ws.Cells.Item(26, 14).Value = 'Squash'
$ws.Cells.Item(26, 15).Value = 'Softball'
$ws.Cells.Item(26, 16).Value = 'Soccer'
$ws.Cells.Item(26, 17).Value = 'Basketball'
$ws.Cells.Item(26, 18).Value = 'Tennis'
$ws.Cells.Item(26, 19).Value = 'Baseball'
$ws.Cells.Item(26, 20).Value = 'Golf'
$ws.Cells.Item(26, 21).Value = 'Running'

